# Update Sema3c-Nrp2 LR-pairs sheet with new TPM-derived values.
# The three "ECs" sending-cluster rows are dropped entirely, and the
# remaining FAPs/MuSCs sending-cluster rows move up to rows 2-7 with
# refreshed numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 2-4 (sending cluster = ECs); everything below
# shifts up to fill the gap.
$ws.Range("A2:T4").EntireRow.Delete() | Out-Null

# New data (rows now at 2-7) in column order A..T
$data = @(
    @("FAPs",  "Sema3c", "Nrp2", "ECs",   3, 1, 42.071953, 126.215859, 0.978774012990499,  0.978774012990499,  3, 1, 25.37147633333333, 76.114429, 0.5780881462719274, 0.5780881462719274, 1067.427559836612, 9606.848038529512, 0.565817654788813,  0.565817654788813),
    @("FAPs",  "Sema3c", "Nrp2", "FAPs",  3, 1, 42.071953, 126.215859, 0.978774012990499,  0.978774012990499,  3, 1, 10.21969166666667, 30.659075, 0.2328552951919536, 0.2328552951919536, 429.9623874744917, 3869.661487270425, 0.2279127117211157, 0.2279127117211157),
    @("FAPs",  "Sema3c", "Nrp2", "MuSCs", 3, 1, 42.071953, 126.215859, 0.978774012990499,  0.978774012990499,  3, 1, 8.297426666666667, 24.89228,  0.189056558536119,  0.189056558536119,  349.0889447409467, 3141.80050266852,  0.1850436464805703, 0.1850436464805704),
    @("MuSCs", "Sema3c", "Nrp2", "ECs",   3, 1, 0.912385,  2.737155,   0.021225987009501,  0.021225987009501,  3, 1, 25.37147633333333, 76.114429, 0.5780881462719274, 0.5780881462719274, 23.14855443438833, 208.336989909495,  0.01227049148311444,0.01227049148311444),
    @("MuSCs", "Sema3c", "Nrp2", "FAPs",  3, 1, 0.912385,  2.737155,   0.021225987009501,  0.021225987009501,  3, 1, 10.21969166666667, 30.659075, 0.2328552951919536, 0.2328552951919536, 9.324293381291668, 83.91864043162501, 0.004942583470837927,0.004942583470837927),
    @("MuSCs", "Sema3c", "Nrp2", "MuSCs", 3, 1, 0.912385,  2.737155,   0.021225987009501,  0.021225987009501,  3, 1, 8.297426666666667, 24.89228,  0.189056558536119,  0.189056558536119,  7.570447629266667, 68.1340286634,     0.004012912055548626,0.004012912055548626)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
